$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Cloe Mercader Pera"
$ws.Range("B2").Value = 13548
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = "24542949V"
$ws.Range("E2").Value = 620108348
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "C. Marita Pujadas 10`nJaén, 42052"
$ws.Range("H2").Value = "hugocanovas@yahoo.com"
$ws.Range("I2").Value = "Elena"
$ws.Range("J2").Value = "terapia"
$ws.Range("K2").Value = 16345
$ws.Range("L2").Value = ""

# Row 3
$ws.Range("A3").Value = "Máximo Agustí Galiano"
$ws.Range("B3").Value = 31127
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = "96746138L"
$ws.Range("E3").Value = 601327677
$ws.Range("F3").Value = 652893194
$ws.Range("G3").Value = "Rambla de Flavia Otero 77`nGranada, 37337"
$ws.Range("H3").Value = "odalysquevedo@hotmail.com"
$ws.Range("I3").Value = "Oriol"
$ws.Range("J3").Value = "terapia"
$ws.Range("K3").Value = 22319
$ws.Range("L3").Value = "Profesional de medicina alternativa"

# Row 4
$ws.Range("A4").Value = "Juliana Amalia Luís Cerro"
$ws.Range("B4").Value = 5639
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "65416954D"
$ws.Range("E4").Value = 673079480
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = "Rambla de Horacio Girón 57 Puerta 6 `nSoria, 07776"
$ws.Range("H4").Value = "jose-angelzabala@hotmail.com"
$ws.Range("I4").Value = "Psicologo3"
$ws.Range("J4").Value = "terapia"
$ws.Range("K4").Value = 18234
$ws.Range("L4").Value = ""

# Row 5
$ws.Range("A5").Value = "Lola Conesa Agustí"
$ws.Range("B5").Value = 7526
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = "92116835H"
$ws.Range("E5").Value = 647530889
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = "Avenida Arcelia Ordóñez 276`nMelilla, 32974"
$ws.Range("H5").Value = "nataliocobos@hotmail.com"
$ws.Range("I5").Value = "Elena"
$ws.Range("J5").Value = "terapia"
$ws.Range("K5").Value = 37511
$ws.Range("L5").Value = ""

# Row 6
$ws.Range("A6").Value = "Serafina Hernandez Blanch"
$ws.Range("B6").Value = 29161
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = "49144449G"
$ws.Range("E6").Value = 603837270
$ws.Range("F6").Value = 697918385
$ws.Range("G6").Value = "Vial de Marcio Casares 166`nLa Rioja, 19332"
$ws.Range("H6").Value = "curro84@gmail.com"
$ws.Range("I6").Value = "Oriol"
$ws.Range("J6").Value = "terapia"
$ws.Range("K6").Value = 22907
$ws.Range("L6").Value = "Técnico en aparatos de diagnóstico y tratamiento médico"

# Row 7
$ws.Range("A7").Value = "Osvaldo Nogués Palacios"
$ws.Range("B7").Value = 27157
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = "33841962S"
$ws.Range("E7").Value = 614952390
$ws.Range("F7").Value = 649549202
$ws.Range("G7").Value = "Calle de Serafina Llanos 4`nBurgos, 46605"
$ws.Range("H7").Value = "vicenteadadia@hotmail.com"
$ws.Range("I7").Value = "Oriol"
$ws.Range("J7").Value = "terapia"
$ws.Range("K7").Value = 2836
$ws.Range("L7").Value = "Compositor"

# Row 8
$ws.Range("A8").Value = "Agustina Vélez Blázquez"
$ws.Range("B8").Value = 35263
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = "80825285L"
$ws.Range("E8").Value = 666341079
$ws.Range("F8").Value = 669667855
$ws.Range("G8").Value = "Camino de Araceli Alegria 88`nSevilla, 04254"
$ws.Range("H8").Value = "kcoronado@yahoo.com"
$ws.Range("I8").Value = "Oriol"
$ws.Range("J8").Value = "terapia"
$ws.Range("K8").Value = 32237
$ws.Range("L8").Value = "Deportista"

# Row 9
$ws.Range("A9").Value = "Cristian Carbajo Hernandez"
$ws.Range("B9").Value = 32731
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = "25820527Z"
$ws.Range("E9").Value = 604199947
$ws.Range("F9").Value = ""
$ws.Range("G9").Value = "Callejón Édgar Calderón 65 Piso 5 `nVizcaya, 10707"
$ws.Range("H9").Value = "anitavergara@yahoo.com"
$ws.Range("I9").Value = "Psicologo3"
$ws.Range("J9").Value = "terapia"
$ws.Range("K9").Value = 22672
$ws.Range("L9").Value = ""

# Row 10
$ws.Range("A10").Value = "Maximiano Canals-Aramburu"
$ws.Range("B10").Value = 12198
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = "15318660Q"
$ws.Range("E10").Value = 670341450
$ws.Range("F10").Value = ""
$ws.Range("G10").Value = "Pasaje de Cornelio Cuesta 3`nSoria, 02873"
$ws.Range("H10").Value = "custodiolaguna@yahoo.com"
$ws.Range("I10").Value = "Elena"
$ws.Range("J10").Value = "terapia"
$ws.Range("K10").Value = 28165
$ws.Range("L10").Value = ""

# Row 11
$ws.Range("A11").Value = "Ciriaco de Acosta"
$ws.Range("B11").Value = 22178
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = "74962378L"
$ws.Range("E11").Value = 651793753
$ws.Range("F11").Value = ""
$ws.Range("G11").Value = "Cañada de Dorita Villegas 51`nBurgos, 03277"
$ws.Range("H11").Value = "joan25@yahoo.com"
$ws.Range("I11").Value = "Oriol"
$ws.Range("J11").Value = "terapia"
$ws.Range("K11").Value = 7173
$ws.Range("L11").Value = ""
